function Set-TextValue {
    param($ws, $cellRef, $val)
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextValue $ws 'D2' '85.445.56'
Set-TextValue $ws 'E2' '  +4.45%  '
Set-TextValue $ws 'D3' '3.225.08'
Set-TextValue $ws 'E3' '  +2.30%  '
Set-TextValue $ws 'E4' '  +0.22%  '
Set-TextValue $ws 'D5' '206.14'
Set-TextValue $ws 'E5' '  -4.91%  '
Set-TextValue $ws 'D6' '615.07'
Set-TextValue $ws 'E6' '  -0.47%  '
Set-TextValue $ws 'D7' '0.347'
Set-TextValue $ws 'E7' '  +21.08%  '
Set-TextValue $ws 'D8' '1.00'
Set-TextValue $ws 'E8' '  +0.10%  '
Set-TextValue $ws 'D9' '0.635'
Set-TextValue $ws 'E9' '  +9.51%  '
Set-TextValue $ws 'D10' '3.218.73'
Set-TextValue $ws 'E10' '  +2.22%  '
Set-TextValue $ws 'D11' '0.561'
Set-TextValue $ws 'E11' '  -5.85%  '
Set-TextValue $ws 'D12' '0.175'
Set-TextValue $ws 'E12' '  +6.20%  '
Set-TextValue $ws 'D13' '0.0000247'
Set-TextValue $ws 'E13' '  -2.80%  '
Set-TextValue $ws 'D14' '3.815.37'
Set-TextValue $ws 'E14' '  +2.13%  '
Set-TextValue $ws 'D15' '32.86'
Set-TextValue $ws 'E15' '  +2.60%  '
Set-TextValue $ws 'D16' '5.17'
Set-TextValue $ws 'E16' '  -2.09%  '
Set-TextValue $ws 'D17' '85.338.24'
Set-TextValue $ws 'E17' '  +4.39%  '
Set-TextValue $ws 'D18' '3.222.58'
Set-TextValue $ws 'E18' '  +2.36%  '
Set-TextValue $ws 'D19' '13.73'
Set-TextValue $ws 'E19' '  -1.42%  '
Set-TextValue $ws 'D20' '2.91'
Set-TextValue $ws 'E20' '  -9.26%  '
Set-TextValue $ws 'D21' '421.27'
Set-TextValue $ws 'E21' '  -2.75%  '
Set-TextValue $ws 'D22' '8.78'
Set-TextValue $ws 'E22' '  -1.13%  '
Set-TextValue $ws 'D23' '5.15'
Set-TextValue $ws 'E23' '  +0.69%  '
Set-TextValue $ws 'D24' '7.15'
Set-TextValue $ws 'E24' '  -1.43%  '
Set-TextValue $ws 'D25' '12.21'
Set-TextValue $ws 'E25' '  +2.76%  '
Set-TextValue $ws 'D26' '5.05'
Set-TextValue $ws 'E26' '  -3.31%  '
Set-TextValue $ws 'D27' '3.395.01'
Set-TextValue $ws 'E27' '  +2.46%  '
Set-TextValue $ws 'D28' '74.78'
Set-TextValue $ws 'E28' '  -2.20%  '
Set-TextValue $ws 'E29' '  +0.18%  '
Set-TextValue $ws 'D30' '0.0000125'
Set-TextValue $ws 'E30' '  +3.36%  '
Set-TextValue $ws 'D31' '0.169'
Set-TextValue $ws 'E31' '  +15.36%  '
Set-TextValue $ws 'E32' '  +0.21%  '
Set-TextValue $ws 'D33' '8.64'
Set-TextValue $ws 'E33' '  -3.78%  '
Set-TextValue $ws 'D34' '535.53'
Set-TextValue $ws 'E34' '  -5.29%  '
Set-TextValue $ws 'D35' '1.40'
Set-TextValue $ws 'E35' '  -6.02%  '
Set-TextValue $ws 'D36' '1.92'
Set-TextValue $ws 'E36' '  -3.07%  '
Set-TextValue $ws 'D37' '6.60'
Set-TextValue $ws 'E37' '  +8.28%  '
Set-TextValue $ws 'D38' '0.134'
Set-TextValue $ws 'E38' '  -11.38%  '
Set-TextValue $ws 'D39' '22.00'
Set-TextValue $ws 'E39' '  -2.40%  '
Set-TextValue $ws 'D40' '0.998'
Set-TextValue $ws 'E40' '  -0.05%  '
Set-TextValue $ws 'D41' '21.60'
Set-TextValue $ws 'D42' '0.384'
Set-TextValue $ws 'E42' '  -4.81%  '
Set-TextValue $ws 'D43' '1.94'
Set-TextValue $ws 'E43' '  -2.97%  '
Set-TextValue $ws 'B44' 'Monero'
Set-TextValue $ws 'C44' 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws 'D44' '157.55'
Set-TextValue $ws 'E44' '  -0.58%  '
Set-TextValue $ws 'B45' 'USDe'
Set-TextValue $ws 'C45' 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue $ws 'D45' '1.00'
Set-TextValue $ws 'E45' '  -0.17%  '
Set-TextValue $ws 'D46' '2.84'
Set-TextValue $ws 'E46' '  -5.80%  '
Set-TextValue $ws 'B47' 'OKB'
Set-TextValue $ws 'C47' 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue $ws 'D47' '43.89'
Set-TextValue $ws 'E47' '  -1.20%  '
Set-TextValue $ws 'B48' 'Aave'
Set-TextValue $ws 'C48' 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws 'D48' '173.41'
Set-TextValue $ws 'E48' '  -6.84%  '
Set-TextValue $ws 'D49' '1.29'
Set-TextValue $ws 'E49' '  -2.25%  '
Set-TextValue $ws 'D50' '4.17'
Set-TextValue $ws 'E50' '  +0.01%  '
Set-TextValue $ws 'D51' '0.721'
Set-TextValue $ws 'E51' '  -5.05%  '
